$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 1. Insert a new column before X (column 24). This shifts every
#    existing column at/after X one letter to the right (Y->Z, Z->AA,
#    AA->AB, AB->AC, AC->AD, AD->AE) and widens the "W" custom-width
#    run to cover the new column too.
# ------------------------------------------------------------------
$ws.Columns("X:X").Insert()
$ws.Columns("X:X").ColumnWidth = $ws.Columns("W:W").ColumnWidth

# ------------------------------------------------------------------
# 2. New "Inhibitor" helper column (X) - header labels + per-row
#    stoichiometry coefficients for the two new reactions.
# ------------------------------------------------------------------
$ws.Range("X1").Value = "Inhibitor"
$ws.Range("X2").Value = "x22"

$ws.Range("X3").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("X5").Value = 0
$ws.Range("X6").Value = 0
$ws.Range("X7").Value = 0
$ws.Range("X8").Value = 0
$ws.Range("X9").Value = 0
$ws.Range("X10").Value = 0
$ws.Range("X11").Value = 0
$ws.Range("X12").Value = 0
$ws.Range("X13").Value = 0
$ws.Range("X14").Value = 0
$ws.Range("X15").Value = 0
$ws.Range("X16").Value = 0
$ws.Range("X17").Value = 0
$ws.Range("X18").Value = 0
$ws.Range("X19").Value = 0
$ws.Range("X20").Value = 0
$ws.Range("X21").Value = 0
$ws.Range("X22").Value = 0
$ws.Range("X23").Value = 0
$ws.Range("X24").Value = 0
$ws.Range("X25").Value = -1
$ws.Range("X26").Value = 1

# ------------------------------------------------------------------
# 3. Other value corrections to existing reaction rows.
# ------------------------------------------------------------------
$ws.Range("C4").Value = 0
$ws.Range("S4").Value = 1

$ws.Range("AB6").Value = 0.1
$ws.Range("AB8").Value = 0.1
$ws.Range("AB10").Value = 0.1
$ws.Range("AB12").Value = 0.1
$ws.Range("AB14").Value = 0.1
$ws.Range("AB15").Value = 0.0001
$ws.Range("AB17").Value = 0.0001

$ws.Range("T22").Value = 1
$ws.Range("U22").Value = 0

# ------------------------------------------------------------------
# 4. Fix mislabeled reaction names in the secondary name/rate table.
# ------------------------------------------------------------------
$ws.Range("Z12").Value = "MEK deactivation"
$ws.Range("Z23").Value = "SOS activation by actTGFa_EGFR"

# ------------------------------------------------------------------
# 5. Brand-new reactions (rows 25 & 26): Inhibition of target (k23)
#    and Increasing Inhib concentration (k24).
# ------------------------------------------------------------------
$ws.Range("A25").Value = "Inhibition of target"
$ws.Range("B25").Value = "k23"
$ws.Range("C25:W25").Value = 0
$ws.Range("D25").Value = -1

$ws.Range("Z25").Value = "Inhibition of target"
$ws.Range("AA25").Value = "k23"
$ws.Range("AB25").Value = 0.01
$ws.Range("AC25").Value = 1

$ws.Range("A26").Value = "Increasing Inhib concentration"
$ws.Range("B26").Value = "k24"
$ws.Range("C26:W26").Value = 0

$ws.Range("Z26").Value = "Increasing Inhib concentration"
$ws.Range("AA26").Value = "k24"
$ws.Range("AB26").Value = 0.5
$ws.Range("AC26").Value = 0

# ------------------------------------------------------------------
# 6. View state: scroll back to the top-left corner and move the
#    selection onto the newly added data.
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J25").Select()
